# Update creation date values in the "Backup Resource Report" sheet.
# The original workbook stores the date strings
#   "2023-12-10T21:00:00Z" and "2023-12-10T22:00:00Z"
# in the shared string table; every cell that references those strings
# must be updated to the new values
#   "2023-12-11T21:00:00Z" and "2023-12-11T22:00:00Z"
# respectively (day incremented from 10 to 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value()

        if ($val -eq "2023-12-10T21:00:00Z") {
            $cell.Value = "2023-12-11T21:00:00Z"
        }
        elseif ($val -eq "2023-12-10T22:00:00Z") {
            $cell.Value = "2023-12-11T22:00:00Z"
        }
    }
}

$wb.Save()
